$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DWH")

# --- Fill in the new CUST_STATS mapping block (rows 51-59) ---
# The write order below intentionally matches the order in which the
# original author entered the values, since that determines the order
# new entries are appended to the shared-string table.

$ws.Cells.Item(52, 2).Value = "price_mean"
$ws.Cells.Item(53, 2).Value = "price_max"
$ws.Cells.Item(54, 2).Value = "price_min"
$ws.Cells.Item(55, 2).Value = "price_median"
$ws.Cells.Item(56, 2).Value = "total_sales"

$ws.Cells.Item(55, 5).Value = "median of all prices bought"
$ws.Cells.Item(53, 5).Value = "maximum of all prices"

$ws.Cells.Item(58, 2).Value = "max_ordersize"

$ws.Cells.Item(54, 5).Value = "minimum of prirces bought"
$ws.Cells.Item(52, 5).Value = "avg of all prices bought"
$ws.Cells.Item(56, 5).Value = "total sales (incl. Taxes)"
$ws.Cells.Item(57, 5).Value = "total number of orders"

$ws.Cells.Item(57, 2).Value = "order count"
$ws.Cells.Item(59, 2).Value = "max_quantity_per_order"

$ws.Cells.Item(51, 1).Value = "CUST_STATS"

$ws.Cells.Item(58, 5).Value = "number of lines in the order with the most lines"
$ws.Cells.Item(59, 5).Value = "maximum total quantity bought in one order"

$ws.Cells.Item(51, 5).Value = "surrogate key of customer"

# Remaining cells (reuse already-known shared strings)
$ws.Cells.Item(51, 2).Value = "CUST_ID"
$ws.Cells.Item(51, 3).Value = "integer"
$ws.Cells.Item(52, 3).Value = "numeric"
$ws.Cells.Item(53, 3).Value = "numeric"
$ws.Cells.Item(54, 3).Value = "numeric"
$ws.Cells.Item(55, 3).Value = "numeric"
$ws.Cells.Item(56, 3).Value = "numeric"
$ws.Cells.Item(57, 3).Value = "numeric"
$ws.Cells.Item(58, 3).Value = "numeric"
$ws.Cells.Item(59, 3).Value = "numeric"

$ws.Cells.Item(52, 1).Value = "CUST_STATS"
$ws.Cells.Item(53, 1).Value = "CUST_STATS"
$ws.Cells.Item(54, 1).Value = "CUST_STATS"
$ws.Cells.Item(55, 1).Value = "CUST_STATS"
$ws.Cells.Item(56, 1).Value = "CUST_STATS"
$ws.Cells.Item(57, 1).Value = "CUST_STATS"
$ws.Cells.Item(58, 1).Value = "CUST_STATS"
$ws.Cells.Item(59, 1).Value = "CUST_STATS"

# --- Column width adjustments ---
# (ColumnWidth inputs chosen so the stored OOXML "width" attribute comes out
# as close as possible to the target 15 / 24.54296875 character widths)
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(2).ColumnWidth = 23.666666666666668

# --- Scroll position / selection state ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E45").Select()
